$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.1946700962556722
$ws.Range("D3").Value = 0.9302171857461253
$ws.Range("E3").Value = 0.2920847147779689
$ws.Range("G3").Value = 'max\_depth: 2, max\_features: 2 \\'
$ws.Range("C4").Value = 0.1744698823805305
$ws.Range("D4").Value = 0.8336919029018337
$ws.Range("E4").Value = 0.270562693368184
$ws.Range("E5").Value = 0.3465710850029271
$ws.Range("C6").Value = 0.2131805773361838
$ws.Range("D6").Value = 1.018668200815774
$ws.Range("E6").Value = 0.3674553873753294
$ws.Range("C7").Value = 0.180004985735156
$ws.Range("D7").Value = 0.8601409999351622
$ws.Range("E7").Value = 0.2606143658648037
$ws.Range("C10").Value = 0.1896974681324275
$ws.Range("D10").Value = 0.9064558365325736
$ws.Range("E10").Value = 0.3110517701337139
$ws.Range("C11").Value = 0.1854587131980305
$ws.Range("D11").Value = 0.8862012480673596
$ws.Range("E11").Value = 0.295901159125759
$ws.Range("C12").Value = 0.2321910471763878
$ws.Range("D12").Value = 1.109508376552066
$ws.Range("E12").Value = 0.2960605863698443
$ws.Range("C13").Value = 0.2401988028572418
$ws.Range("D13").Value = 1.147772866562917
$ws.Range("E13").Value = 0.2951683438880019
$ws.Range("C14").Value = 0.2350197726360835
$ws.Range("D14").Value = 1.123025239629539
$ws.Range("E14").Value = 0.2955565190977457
$ws.Range("C15").Value = 0.1717752899257359
$ws.Range("D15").Value = 0.8208159848320152
$ws.Range("E15").Value = 0.2642543000809127
$ws.Range("C16").Value = 0.2181644356426673
$ws.Range("D16").Value = 1.042483212659841
$ws.Range("E16").Value = 0.3804028698510015
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 'omega: 1000, r: 0.1 \\'
$ws.Range("C17").Value = 0.2181644356426673
$ws.Range("D17").Value = 1.042483212659841
$ws.Range("E17").Value = 0.3804028698510015
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 'omega: 1000, r: 0.1 \\'
$ws.Range("C18").Value = 0.2083008123370609
$ws.Range("D18").Value = 0.9953505914248333
$ws.Range("E18").Value = 0.3552823554458325
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 'mu: 0.5, omega: 1000 \\'
$ws.Range("C19").Value = 0.226124091023563
$ws.Range("D19").Value = 1.080517858814232
$ws.Range("E19").Value = 0.3978103973833274
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 'alpha: 0.001, beta: 0.1, lambda1: 0.001, omega: 100, sigma: 0.1 \\'
$ws.Range("C20").Value = 0.2556219654169771
$ws.Range("D20").Value = 1.221471350036101
$ws.Range("E20").Value = 0.4231730705141653
$ws.Range("F20").Value = 390
$ws.Range("G20").Value = 'alpha: 0.01, lambda1: 0.5, omega: 10000, sigma: 0.001, w: 50 \\'
$ws.Range("C21").Value = 0.2275338911951721
$ws.Range("D21").Value = 1.087254488493483
$ws.Range("E21").Value = 0.4001770000827641
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 'alpha: 0.001, beta: 0.01, e\_utility: 0.05, lambda1: 0.75, omega: 10000, pi: 0.5, sigma: 0.25 \\'
$ws.Range("C22").Value = 0.1884916474986744
$ws.Range("D22").Value = 0.9006939085425074
$ws.Range("E22").Value = 0.2855417896556091
$ws.Range("F22").Value = 19
$ws.Range("G22").Value = 'alpha: 0.1, beta: 0.25, e\_utility: 0.05, lambda1: 0.001, sigma: 50 \\'
$ws.Range("C23").Value = 0.2506050574535119
$ws.Range("D23").Value = 1.19749841276076
$ws.Range("E23").Value = 0.3004655666067801
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 'fuzzy\_operator: min, rules: 10 \\'
$ws.Range("C24").Value = 0.1917625040853879
$ws.Range("D24").Value = 0.9163234637109388
$ws.Range("E24").Value = 0.3169204694247891
$ws.Range("C25").Value = 0.2219726452010427
$ws.Range("D25").Value = 1.060680470719809
$ws.Range("E25").Value = 0.3885312758881269
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 'adaptive\_filter: wRLS, fuzzy\_operator: prod, rules: 1 \\'
$ws.Range("C26").Value = 0.1919697938121079
$ws.Range("D26").Value = 0.9173139829017769
$ws.Range("E26").Value = 0.2707153823630994
$ws.Range("F26").Value = 17
$ws.Range("G26").Value = 'error\_metric: RMSE, fuzzy\_operator: prod, num\_generations: 10, num\_parents\_mating: 5, parallel\_processing: 10, rules: 17, sol\_per\_pop: 10 \\'
$ws.Range("C27").Value = 0.1920445580748718
$ws.Range("D27").Value = 0.9176712386048368
$ws.Range("E27").Value = 0.290631504079077
$ws.Range("G27").Value = 'adaptive\_filter: RLS, error\_metric: MAE, fuzzy\_operator: max, lambda1: 0.98, num\_generations: 5, num\_parents\_mating: 5, parallel\_processing: 10, rules: 1, sol\_per\_pop: 5 \\'
$ws.Range("C28").Value = 0.1806897359566609
$ws.Range("D28").Value = 0.863413030083803
$ws.Range("E28").Value = 0.2563998979109234
$ws.Range("G28").Value = 'adaptive\_filter: wRLS, error\_metric: MAE, fuzzy\_operator: max, num\_generations: 5, num\_parents\_mating: 5, parallel\_processing: 10, rules: 19, sol\_per\_pop: 5 \\'
$ws.Range("C29").Value = 0.2101059302043536
$ws.Range("D29").Value = 1.00397621854862
$ws.Range("E29").Value = 0.2800609286486028
$ws.Range("G29").Value = 'combination: weighted\_average, n\_estimators: 50 \\'
$ws.Range("C30").Value = 0.1777690567093121
$ws.Range("D30").Value = 0.8494567723832487
$ws.Range("E30").Value = 0.2582718002045762
$ws.Range("G30").Value = 'combination: median, n\_estimators: 50 \\'
$ws.Range("C31").Value = 0.1731568207039562
$ws.Range("D31").Value = 0.827417531229002
$ws.Range("E31").Value = 0.2602497450332074
